# Weekly Fruta/hortaliza update for Nectarin sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the date-format used by column D (style s="2") so new/rewritten
# date cells keep the same number format as the rest of the column.
$dateFormat = $ws.Range("D284").NumberFormat

# --- Update existing rows 284-286 with the new (current week) records ---

# Row 284: Red Diamond / Especial
$ws.Cells.Item(284, 4).Value = 44595
$ws.Cells.Item(284, 4).NumberFormat = $dateFormat
$ws.Cells.Item(284, 11).Value = "Red Diamond"
$ws.Cells.Item(284, 12).Value = "Especial"
$ws.Cells.Item(284, 13).Value = 100
$ws.Cells.Item(284, 15).Value = 13000
$ws.Cells.Item(284, 16).Value = 13000
$ws.Cells.Item(284, 19).Value = 812

# Row 285: Red Diamond / Primera
$ws.Cells.Item(285, 4).Value = 44595
$ws.Cells.Item(285, 4).NumberFormat = $dateFormat
$ws.Cells.Item(285, 11).Value = "Red Diamond"
$ws.Cells.Item(285, 14).Value = 11000
$ws.Cells.Item(285, 15).Value = 12000
$ws.Cells.Item(285, 16).Value = 11500
$ws.Cells.Item(285, 19).Value = 719

# Row 286: Venus / Especial
$ws.Cells.Item(286, 4).Value = 44595
$ws.Cells.Item(286, 4).NumberFormat = $dateFormat
$ws.Cells.Item(286, 11).Value = "Venus"
$ws.Cells.Item(286, 12).Value = "Especial"
$ws.Cells.Item(286, 13).Value = 100

# --- Append new rows 287-290 ---
# Row 287: Venus / Primera (new)
# Row 288: Super Queen / Primera (old row 284 data, date 44552)
# Row 289: Super Queen / Primera (old row 285 data, date 44544)
# Row 290: Super Queen / Segunda (old row 286 data, date 44544)

$newRows = @(
    @{ Row = 287; D = 44595; K = "Venus";       L = "Primera"; M = 120; N = 11000; O = 12000; P = 11500; S = 719 },
    @{ Row = 288; D = 44552; K = "Super Queen"; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; S = 844 },
    @{ Row = 289; D = 44544; K = "Super Queen"; L = "Primera"; M = 120; N = 14000; O = 15000; P = 14500; S = 906 },
    @{ Row = 290; D = 44544; K = "Super Queen"; L = "Segunda"; M = 60;  N = 13000; O = 13000; P = 13000; S = 812 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value = "Ñuble"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 5).Value = 16
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103006
    $ws.Cells.Item($row, 10).Value = "Nectarín"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
